$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" column (C) values for rows 2-252 based on the
# new segmented pattern described in the commit diff.
# Rows 2-19   (Generation 0-17):  7343
# Rows 20-61  (Generation 18-59): 7310
# Rows 62-252 (Generation 60-250): 7293

$ws.Range("C2:C19").Value = 7343
$ws.Range("C20:C61").Value = 7310
$ws.Range("C62:C252").Value = 7293
